$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7 (columns A:T), reflecting updated TPM values and
# two additional rows (5 and 7) inserted for the "FAPs" target cluster.
$data = @(
    @("ECs",   "Gast", "Cckbr", "ECs",  3, 1,                  3.158628666666667,  9.475885999999999,  0.6322622710965859, 0.6322622710965859, 1, 0.3333333333333333, 0.09168033333333332, 0.275041,  0.1820807603088688, 0.1820807603088689, 0.2895841290362222,  2.606257161325999,  0.1151227950358785,  0.1151227950358785),
    @("ECs",   "Gast", "Cckbr", "FAPs", 3, 1,                  3.158628666666667,  9.475885999999999,  0.6322622710965859, 0.6322622710965859, 3, 1,                  0.4118343333333334,  1.235503,  0.8179192396911312, 0.8179192396911312, 1.300831731184222,   11.707485580658,    0.5171394760607074,  0.5171394760607074),
    @("FAPs",  "Gast", "Cckbr", "ECs",  1, 0.3333333333333333, 1.206677666666667,  3.620033,            0.2415405046055416, 0.2415405046055416, 1, 0.3333333333333333, 0.09168033333333332, 0.275041,  0.1820807603088688, 0.1820807603088689, 0.1106286107058889,  0.9956574963529998, 0.04397987872396485, 0.04397987872396485),
    @("FAPs",  "Gast", "Cckbr", "FAPs", 1, 0.3333333333333333, 1.206677666666667,  3.620033,            0.2415405046055416, 0.2415405046055416, 3, 1,                  0.4118343333333334,  1.235503,  0.8179192396911312, 0.8179192396911312, 0.4969512923998889,  4.472561631599,     0.1975606258815768,  0.1975606258815768),
    @("MuSCs", "Gast", "Cckbr", "ECs",  2, 0.6666666666666666, 0.6304506666666666, 1.891352,            0.1261972242978725, 0.1261972242978725, 1, 0.3333333333333333, 0.09168033333333332, 0.275041,  0.1820807603088688, 0.1820807603088689, 0.05779992727022221, 0.520199345432,     0.02297808654902548, 0.02297808654902549),
    @("MuSCs", "Gast", "Cckbr", "FAPs", 2, 0.6666666666666666, 0.6304506666666666, 1.891352,            0.1261972242978725, 0.1261972242978725, 3, 1,                  0.4118343333333334,  1.235503,  0.8179192396911312, 0.8179192396911312, 0.2596412300062222,  2.336771070056,     0.103219137748847,   0.1032191377488471)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

$wb.Save()
